$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26, shifting existing rows 26:103 down to 27:104
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new week's data
$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value = "Ñuble"
$ws.Cells.Item(26, 4).Value = 44672
$ws.Cells.Item(26, 5).Value = 16
$ws.Cells.Item(26, 6).Value = 100112030
$ws.Cells.Item(26, 7).Value = "Poroto granado"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 60
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 21000
$ws.Cells.Item(26, 13).Value = 20500
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Región del Maule"
$ws.Cells.Item(26, 16).Value = 820
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
